# Actualización desde MV -datos-
# Append 30 new daily rows (05-08-2021 .. 03-09-2021) to Sheet1, mirroring the
# existing pattern of columns B/C/D.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$dates = @(
    "05-08-2021","06-08-2021","07-08-2021","08-08-2021","09-08-2021",
    "10-08-2021","11-08-2021","12-08-2021","13-08-2021","14-08-2021",
    "15-08-2021","16-08-2021","17-08-2021","18-08-2021","19-08-2021",
    "20-08-2021","21-08-2021","22-08-2021","23-08-2021","24-08-2021",
    "25-08-2021","26-08-2021","27-08-2021","28-08-2021","29-08-2021",
    "30-08-2021","31-08-2021","01-09-2021","02-09-2021","03-09-2021"
)

$startRow = 218
$endRow = $startRow + $dates.Length - 1
$colARange = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 1))

# Force column A to text before typing so ambiguous day/month date strings
# (e.g. "05-08-2021") are kept as literal text instead of being parsed into
# date serials, matching the source data (plain text dates).
$colARange.NumberFormat = "@"

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = 17537
    if ($i -lt 22) {
        $ws.Cells.Item($row, 3).Value = 2392
    } else {
        $ws.Cells.Item($row, 3).Value = 1456
    }
    $ws.Cells.Item($row, 4).Value = 521
}

# Drop the temporary text number-format so these cells end up with the same
# (default/no explicit style) formatting as the rest of column A.
$colARange.ClearFormats()
